$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for new columns I (I0) and J (IF), rows 2-23
$values = @(
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 9),
    @(10, 10),
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(4, 4),
    @(4, 4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
